# Loan RBI, Variable Instalments
# The "Repayment Schedule" sheet gains a new blank column between the
# existing "In Advance" (M) and "Late" (N) columns, pushing "Late" from N
# to O and "Outstanding" from P to Q. The new column keeps the same
# width as its neighbour "In Advance" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at N - this shifts the old N ("Late") to O
# and the old P ("Outstanding") to Q, exactly matching the target layout.
$ws.Columns("N").Insert()

# Give the freshly inserted column N a fixed custom width (matches the
# "Paid" / "In Advance" style columns near it).
$ws.Columns("N").ColumnWidth = 10.307291666666666

# Restore/update the active selection on the sheet.
$ws.Range("T10").Select()
